# Update the player roster table (A2:C18) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Anfernee Simons",        "PG,SG",       "Portland Trail Blazers"),
    @("James Harden",           "PG,SG",       "LA Clippers"),
    @("Donte DiVincenzo",       "PG,SG,SF",    "Minnesota Timberwolves"),
    @("Keyonte George",         "PG,SG",       "Utah Jazz"),
    @("Quentin Grimes",         "SG,SF",       "Philadelphia 76ers"),
    @("Amen Thompson",          "PG,SG,SF,PF", "Houston Rockets"),
    @("RJ Barrett",             "SG,SF,PF",    "Toronto Raptors"),
    @("Bradley Beal",           "PG,SG,SF",    "Phoenix Suns"),
    @("Draymond Green",         "PF,C",        "Golden State Warriors"),
    @("Jaren Jackson Jr.",      "PF,C",        "Memphis Grizzlies"),
    @("Jayson Tatum",           "SF,PF",       "Boston Celtics"),
    @("Ivica Zubac",            "C",           "LA Clippers"),
    @("Anthony Edwards",        "SG,SF",       "Minnesota Timberwolves"),
    @("Zion Williamson",        "PF,C",        "New Orleans Pelicans"),
    @("Giannis Antetokounmpo",  "PF,C",        "Milwaukee Bucks"),
    @("Fred VanVleet",          "PG",          "Houston Rockets"),
    @("Paul George",            "SG,SF,PF",    "Philadelphia 76ers")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
